# Daily attendance processing - 2026-01-01 22:32:23
# Normalises the "Recorded By" column (G) so the most-recent recorder name
# is listed first. This re-orders the comma-separated list of recorder
# identities for a fixed set of known combinations found in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# Exact "old combination" -> "new combination" replacements for the
# "Recorded By" column (column G).
$map = @{
    "backup@backdoor.com, System, system" = "system, backup@backdoor.com, System"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($null -eq $current) { continue }

    if ($map.ContainsKey($current)) {
        $cell.Value2 = $map[$current]
    }
}
